# "aggiornamento fino a 28 luglio" - append the new daily rows (302-328)
# below the existing data in Sheet1 (date serial, nuovi pos., somma mobile
# 7gg., somma mobile 7gg. per 100mila abitanti), matching the style/format
# already used by the previous rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row holding the last existing record (A:D) - used as the formatting
# template (date style, border, bold, etc.) for the newly appended rows.
$lastRow = 301
$templateDateCell = $ws.Cells.Item($lastRow, 1)

# New daily data: day, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$newData = @(
    @(44376, 1, 2, 18.42468908337172),
    @(44377, 0, 2, 18.42468908337172),
    @(44378, 1, 3, 27.63703362505758),
    @(44379, 0, 3, 27.63703362505758),
    @(44380, 1, 4, 36.84937816674343),
    @(44381, 0, 4, 36.84937816674343),
    @(44382, 0, 3, 27.63703362505758),
    @(44383, 0, 2, 18.42468908337172),
    @(44384, 0, 2, 18.42468908337172),
    @(44385, 0, 1, 9.212344541685859),
    @(44386, 0, 1, 9.212344541685859),
    @(44387, 0, 0, 0),
    @(44388, 0, 0, 0),
    @(44389, 0, 0, 0),
    @(44390, 0, 0, 0),
    @(44391, 0, 0, 0),
    @(44392, 0, 0, 0),
    @(44393, 1, 1, 9.212344541685859),
    @(44394, 0, 1, 9.212344541685859),
    @(44395, 0, 1, 9.212344541685859),
    @(44396, 2, 3, 27.63703362505758),
    @(44397, 0, 3, 27.63703362505758),
    @(44398, 0, 3, 27.63703362505758),
    @(44399, 0, 3, 27.63703362505758),
    @(44400, 1, 3, 27.63703362505758),
    @(44401, 0, 3, 27.63703362505758),
    @(44402, 0, 3, 27.63703362505758)
)

$row = $lastRow + 1
foreach ($rec in $newData) {
    $dateCell = $ws.Cells.Item($row, 1)

    # Copy the date column's style (border, bold font, center/top alignment,
    # custom date number format) from the previous row, then overwrite the
    # value so the applied formatting is kept rather than the source value.
    $templateDateCell.Copy($dateCell)
    $dateCell.Value = $rec[0]

    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]

    $row = $row + 1
}

Write-Host "Appended rows $($lastRow + 1)-$($row - 1) (aggiornamento fino a 28 luglio)"
